$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 19 is new; copy the date-column formatting (style) from A18 before overwriting values
$ws.Range("A18").Copy()
$ws.Range("A19").PasteSpecial(-4122)

$ws.Range("A2").Value = 39400
$ws.Range("B2").Value = 2007
$ws.Range("C2").Value = 0.4235526809466261
$ws.Range("D2").Value = 2008
$ws.Range("E2").Value = 0.6439341879002525

$ws.Range("A3").Value = 39765
$ws.Range("B3").Value = 2008
$ws.Range("C3").Value = -0.5718076928962645
$ws.Range("D3").Value = 2009
$ws.Range("E3").Value = -0.1800933741311961

$ws.Range("A4").Value = 40130
$ws.Range("B4").Value = 2009
$ws.Range("C4").Value = 0.3486139762225005
$ws.Range("D4").Value = 2010
$ws.Range("E4").Value = 0.1555182634501051

$ws.Range("A5").Value = 40494
$ws.Range("B5").Value = 2010
$ws.Range("C5").Value = -0.1384957661262898
$ws.Range("D5").Value = 2011
$ws.Range("E5").Value = 0.6938817570587785

$ws.Range("A6").Value = 40862
$ws.Range("B6").Value = 2011
$ws.Range("C6").Value = 1.566479473280147
$ws.Range("D6").Value = 2012
$ws.Range("E6").Value = 0.9614071719361794

$ws.Range("A7").Value = 41228
$ws.Range("B7").Value = 2012
$ws.Range("C7").Value = 0.7307568962936939
$ws.Range("D7").Value = 2013
$ws.Range("E7").Value = 1.09290550768979

$ws.Range("A8").Value = 41592
$ws.Range("B8").Value = 2013
$ws.Range("C8").Value = 0.818818812164257
$ws.Range("D8").Value = 2014
$ws.Range("E8").Value = 0.9607602172681418

$ws.Range("A9").Value = 41957
$ws.Range("B9").Value = 2014
$ws.Range("C9").Value = 0.9180054319587239
$ws.Range("D9").Value = 2015
$ws.Range("E9").Value = 1.375398114243209

$ws.Range("A10").Value = 42321
$ws.Range("B10").Value = 2015
$ws.Range("C10").Value = 1.984684278296656
$ws.Range("D10").Value = 2016
$ws.Range("E10").Value = 1.473274087935805

$ws.Range("A11").Value = 42689
$ws.Range("B11").Value = 2016
$ws.Range("C11").Value = 1.755995812646982
$ws.Range("D11").Value = 2017
$ws.Range("E11").Value = 1.681032827388362

$ws.Range("A12").Value = 43053
$ws.Range("B12").Value = 2017
$ws.Range("C12").Value = 1.946965557828384
$ws.Range("D12").Value = 2018
$ws.Range("E12").Value = 1.755491062323111

$ws.Range("A13").Value = 43418
$ws.Range("B13").Value = 2018
$ws.Range("C13").Value = 1.06432145354225
$ws.Range("D13").Value = 2019
$ws.Range("E13").Value = 0.776718238020746

$ws.Range("A14").Value = 43783
$ws.Range("B14").Value = 2019
$ws.Range("C14").Value = 1.361817904277696
$ws.Range("D14").Value = 2020
$ws.Range("E14").Value = 1.316199564471554

$ws.Range("A15").Value = 44159
$ws.Range("B15").Value = 2020
$ws.Range("C15").Value = -4.352425014431304
$ws.Range("D15").Value = 2021
$ws.Range("E15").Value = 0.03547044462246518

$ws.Range("A16").Value = 44525
$ws.Range("B16").Value = 2021
$ws.Range("C16").Value = -1.761645650979182
$ws.Range("D16").Value = 2022
$ws.Range("E16").Value = 3.765721202592909

$ws.Range("A17").Value = 44890
$ws.Range("B17").Value = 2022
$ws.Range("C17").Value = 5.20787683103745
$ws.Range("D17").Value = 2023
$ws.Range("E17").Value = 3.217995704408838

$ws.Range("A18").Value = 45254
$ws.Range("B18").Value = 2023
$ws.Range("C18").Value = -0.9008525709169546
$ws.Range("D18").Value = 2024
$ws.Range("E18").Value = 0.6027009207580036

$ws.Range("A19").Value = 45618
$ws.Range("B19").Value = 2024
$ws.Range("C19").Value = 0.2738544794132824
$ws.Range("D19").Value = 2025
$ws.Range("E19").Value = 0.2681899963140832

